$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) to reflect the new "through" date
$ws.Name = "Through 2022-03-27"

# Update the header cell text (column I, row 1) which shows "2022 (through 03-26)"
$ws.Range("I1").Value = "2022 (through 03-27)"

# Update March 2022 total (row 4, column I)
$ws.Range("I4").Value = 111

# Update grand Total row (row 14, column I)
$ws.Range("I14").Value = 411
